# notebooks/terrain_data.xlsx — "adding races, removing others"
#
# The canonical-OOXML diff shows a large amount of <v>N</v> shared-string
# index churn across the sheet, but nearly all of it is a side effect of
# one content change: the shared string "Kuurne-Bruxelles-Kuurne" (the
# race name used by rows 95-101, column C) is renamed to
# "Kuurne-Brussel-Kuurne". Deleting the old string and appending the new
# one at the end of sst shifts every higher shared-string index down by
# one, which is exactly the pattern of every other C/D-column <v> change
# in the diff (146->145, 147->146, ... 149->148, etc). That re-indexing
# is handled automatically by the host when the workbook is saved, so it
# does not need to be (and should not be) written out explicitly here —
# only the actual content change does.
#
# The diff's second content change is that rows 95-101 gain explicit 0
# values in the previously-blank terrain-breakdown columns (G:V).
#
# Finally, the diff shows the saved view state moved (scrolled down, and
# the active selection changed from D69 to C22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename the race name string for the Kuurne-Bruxelles-Kuurne rows ---
$kbkRows = 95..101
foreach ($r in $kbkRows) {
    $ws.Cells.Item($r, 3).Value = "Kuurne-Brussel-Kuurne"   # column C
}

# --- Fill newly-populated terrain columns (previously blank) with 0 ---
# Column indices: G=7 ... V=22
$newZeroCells = @{}
$newZeroCells[95]  = 7..22                     # G95:V95
$newZeroCells[96]  = 7..22                     # G96:V96
$newZeroCells[97]  = 11,13,17,19,20,21,22      # K97,M97,Q97,S97:V97
$newZeroCells[98]  = 7..21                     # G98:U98
$newZeroCells[99]  = 7..21                     # G99:U99
$newZeroCells[100] = 7..22                     # G100:V100
$newZeroCells[101] = 11,13,17,19,20,21,22      # K101,M101,Q101,S101:V101

foreach ($r in $newZeroCells.Keys) {
    foreach ($c in $newZeroCells[$r]) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# --- Restore the view state (scroll position / active selection) ---
# TopLeftCell itself isn't writable through this host's Window object,
# but ScrollRow/ScrollColumn are attempted for completeness, and the
# active selection (which IS persisted) is set to match.
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()
